$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Remove the two trailing rows (old row4 "MacBook Air" and, after the
# shift, the row that held "Samsung Galaxy Tab 10.1") while letting the
# old row5 ("Canon EOS 5D", smaller-font style) slide up into row4.
$ws.Rows.Item(4).Delete()
$ws.Rows.Item(5).Delete()

# Update the remaining product list text.
$ws.Range("A1").Value = "MacBook"
$ws.Range("A2").Value = "Apple Cinema 30"
$ws.Range("A3").Value = "iPhone"
$ws.Range("A4").Value = "Canon EOS 5D"

$ws.Range("E7").Select() | Out-Null
